$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 25 new columns before column P (old P:V block shifts to AO:AU)
$ws.Range("P1:AN1").EntireColumn.Insert()

# Fill in the headers for the newly inserted columns (P1:AN1), continuing the
# "Unnamed: 0.1.1..." naming sequence that already runs through B1:O1.
# (The column insert above already carried the bold/border/centered header
# style into these cells, so only the text needs to be written.)
$dots = ".1.1.1.1.1.1.1.1.1.1.1.1.1.1"
for ($i = 0; $i -lt 25; $i++) {
    $col = 16 + $i
    $ws.Cells.Item(1, $col).Value = "Unnamed: 0" + $dots
    $dots = $dots + ".1"
}

# Fill in the data rows for the newly inserted columns, mirroring the
# 0 / 1 pattern used by columns B:O.
for ($i = 0; $i -lt 25; $i++) {
    $col = 16 + $i
    $ws.Cells.Item(2, $col).Value = 0
    $ws.Cells.Item(3, $col).Value = 1
}

# Update the values that changed on row 2 ("Preço Atual" / "Local"),
# now located at AT2 / AU2 after the column insert.
$ws.Range("AT2").Value = 3905.07
$ws.Range("AU2").Value = "Magalu"
